$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericTextCells = @("D4","D5","D6","D7","D9","D10","D11","D12","D13","D14","D16","D21","D24","D25","D26","D28","D29","D30","D31","D32","D33","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D47","D48","D49","D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.244.31'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.904.83'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("D4").Value = '0.9997'
$ws.Range("D5").Value = '306.36'
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").Value = '0.9994'
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").Value = '0.5442'
$ws.Range("E7").Value = '  +4.13%  '
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("D9").Value = '0.07295'
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("D10").Value = '22.19'
$ws.Range("E10").Value = '  +5.30%  '
$ws.Range("D11").Value = '0.9039'
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").Value = '0.08188'
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("D13").Value = '95.80'
$ws.Range("E13").Value = '  -0.23%  '
$ws.Range("D14").Value = '5.351'
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").Value = '14.88'
$ws.Range("E16").Value = '  +2.08%  '
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("D19").Value = '27.273.44'
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("D20").Value = '1.172.18'
$ws.Range("E20").Value = '  -38.58%  '
$ws.Range("D21").Value = '5.050'
$ws.Range("E21").Value = '  -0.60%  '
$ws.Range("E22").Value = '  +1.37%  '
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").Value = '148.58'
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D25").Value = '2.311'
$ws.Range("E25").Value = '  +1.24%  '
$ws.Range("D26").Value = '18.39'
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("D28").Value = '117.00'
$ws.Range("E28").Value = '  +1.75%  '
$ws.Range("D29").Value = '4.859'
$ws.Range("E29").Value = '  +1.56%  '
$ws.Range("D30").Value = '4.691'
$ws.Range("E30").Value = '  -3.23%  '
$ws.Range("D31").Value = '0.09235'
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").Value = '0.8327'
$ws.Range("E32").Value = '  +5.77%  '
$ws.Range("D33").Value = '0.05084'
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("E34").Value = '  +0.91%  '
$ws.Range("D35").Value = '3.009'
$ws.Range("E35").Value = '  +1.76%  '
$ws.Range("D36").Value = '3.322'
$ws.Range("E36").Value = '  -2.95%  '
$ws.Range("D37").Value = '2.696'
$ws.Range("E37").Value = '  +4.25%  '
$ws.Range("D38").Value = '0.5962'
$ws.Range("E38").Value = '  +4.86%  '
$ws.Range("D39").Value = '0.02003'
$ws.Range("E39").Value = '  +1.30%  '
$ws.Range("D40").Value = '1.080'
$ws.Range("E40").Value = '  +0.44%  '
$ws.Range("D41").Value = '9.285'
$ws.Range("E41").Value = '  +2.98%  '
$ws.Range("D42").Value = '6.676'
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("D43").Value = '116.57'
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").Value = '0.5127'
$ws.Range("E44").Value = '  +5.69%  '
$ws.Range("D45").Value = '0.1532'
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("E46").Value = '  +1.55%  '
$ws.Range("D47").Value = '0.9987'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("D48").Value = '1.645'
$ws.Range("E48").Value = '  +1.66%  '
$ws.Range("D49").Value = '38.33'
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("E50").Value = '  +2.81%  '
$ws.Range("D51").Value = '63.59'
$ws.Range("E51").Value = '  +0.19%  '

foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}
